$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "last_edited_time" shared string that Tháng 8 (row 6) now shares with
# Tháng 7..Tháng 2 (rows 7-12) gets updated to the new edit timestamp.
$newEditedTime = "2024-08-03T03:17:00.000Z"
$ws.Range("D6").Value = $newEditedTime
$ws.Range("D7").Value = $newEditedTime
$ws.Range("D8").Value = $newEditedTime
$ws.Range("D9").Value = $newEditedTime
$ws.Range("D10").Value = $newEditedTime
$ws.Range("D11").Value = $newEditedTime
$ws.Range("D12").Value = $newEditedTime

# last_edited_by.id switches to the other user, matching rows 7-12.
$ws.Range("N6").Value = "41cabcaf-915d-46a5-8eff-38727be27269"

# Report figures for Tháng 8 (row 6) get populated, same as the other
# already-filled months.
$ws.Range("T6").Value = 3000000
$ws.Range("W6").Value = 9400000
$ws.Range("AA6").Value = 15700000
$ws.Range("AE6").Value = 25100000
$ws.Range("AH6").Value = 25100000
$ws.Range("AK6").Value = 6
$ws.Range("AN6").Value = 0
$ws.Range("AQ6").Value = 28100000
